$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"
$ws.Range("A2").Value = "Dami Sanyaolu"

$ws.Columns.Item(1).ColumnWidth = 15.67

$ws.Range("A3").Select()
